# Update the configured module-config path from "Configure\Modules.xlsx"
# to "Design\Configure\Modules.xlsx" (the project folder was renamed/moved
# into a "Design" subfolder), and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 holds the path to the Modules.xlsx config file - update it to include
# the new "Design" folder segment.
$ws.Range("B2").Value = "Design\Configure\Modules.xlsx"

# Move the selected/active cell to E4 (reflects where the author left the
# cursor when they saved the file).
[void]$ws.Range("E4").Select()
